$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "No of Cows" / "Area of Field" rows up by one (row5->row4, row6->row5) ---
$ws.Range("C4").Value = "No of Cows"
$ws.Range("D4").Value = 10

$ws.Range("C5").Value = "Area of Field"
$ws.Range("D5").Value = 10

# Remove the now-vacated old row 6 content
$ws.Range("C6:D6").ClearContents()

# --- Biomass formula now references the shifted rows ---
$ws.Range("D10").Formula = "=D4+D5"

# D11 (Energy) formula is unchanged: =D10 *5

# --- New rows further down the sheet ---
$ws.Range("D74").Formula = "=D4+D5"
$ws.Range("D50").Formula = "=D74*5"

# --- Update the view state to match the saved selection/scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("D75").Select()
